$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date style (numFmtId 14) from an existing date cell so we don't
# introduce a brand-new style entry, then set the values/content.
$ws.Cells.Item(18, 1).Copy()
$ws.Cells.Item(19, 1).PasteSpecial(-4122)
$ws.Cells.Item(18, 1).Copy()
$ws.Cells.Item(20, 1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 19: new titration data point
$ws.Cells.Item(19, 1).Value = 43208
$ws.Cells.Item(19, 2).Value = 2224.6526515354399
$ws.Cells.Item(19, 6).Value = "With Junk"

# Row 20: new titration data point
$ws.Cells.Item(20, 1).Value = 43208
$ws.Cells.Item(20, 2).Value = 2223.9979788785199
$ws.Cells.Item(20, 6).Value = "end of sample"

# Update selection to reflect the new active cell
$ws.Range("F21").Select()
